$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width for BK (col 63) = 12, matching existing column style
$ws.Range("BK1").ColumnWidth = 11.14

# BK1 header date label (stored as text, matching style of BJ1)
$ws.Range("BK1").NumberFormat = "@"
$ws.Range("BK1").Value = "2024/11/10"
$ws.Range("BJ1").Copy()
$ws.Range("BK1").PasteSpecial(-4122)

$ws.Range("BK2").Value = 149.6
$ws.Range("A2").Copy()
$ws.Range("BK2").PasteSpecial(-4122)
$ws.Range("BK3").Value = 121.8
$ws.Range("AD3").Copy()
$ws.Range("BK3").PasteSpecial(-4122)
$ws.Range("BK4").Value = 171.1
$ws.Range("A4").Copy()
$ws.Range("BK4").PasteSpecial(-4122)
$ws.Range("BK5").Value = 143.1
$ws.Range("A5").Copy()
$ws.Range("BK5").PasteSpecial(-4122)
$ws.Range("BK6").Value = 145.1
$ws.Range("A6").Copy()
$ws.Range("BK6").PasteSpecial(-4122)
$ws.Range("BK7").Value = 127.3
$ws.Range("B7").Copy()
$ws.Range("BK7").PasteSpecial(-4122)
$ws.Range("BK8").Value = 113.8
$ws.Range("H8").Copy()
$ws.Range("BK8").PasteSpecial(-4122)
$ws.Range("BK9").Value = 154.3
$ws.Range("A9").Copy()
$ws.Range("BK9").PasteSpecial(-4122)
$ws.Range("BK10").Value = 148.3
$ws.Range("A10").Copy()
$ws.Range("BK10").PasteSpecial(-4122)
$ws.Range("BK11").Value = 145.2
$ws.Range("A11").Copy()
$ws.Range("BK11").PasteSpecial(-4122)
$ws.Range("BK12").Value = 115
$ws.Range("C12").Copy()
$ws.Range("BK12").PasteSpecial(-4122)
$ws.Range("BK13").Value = 136.3
$ws.Range("D13").Copy()
$ws.Range("BK13").PasteSpecial(-4122)
$ws.Range("BK14").Value = 174.2
$ws.Range("A14").Copy()
$ws.Range("BK14").PasteSpecial(-4122)
$ws.Range("BK15").Value = 114.9
$ws.Range("P15").Copy()
$ws.Range("BK15").PasteSpecial(-4122)
$ws.Range("BK16").Value = 152.1
$ws.Range("A16").Copy()
$ws.Range("BK16").PasteSpecial(-4122)
$ws.Range("BK17").Value = 135.8
$ws.Range("C17").Copy()
$ws.Range("BK17").PasteSpecial(-4122)
$ws.Range("BK18").Value = 146.3
$ws.Range("A18").Copy()
$ws.Range("BK18").PasteSpecial(-4122)
$ws.Range("BK19").Value = 152.6
$ws.Range("A19").Copy()
$ws.Range("BK19").PasteSpecial(-4122)
$ws.Range("BK20").Value = 231.5
$ws.Range("A20").Copy()
$ws.Range("BK20").PasteSpecial(-4122)
$ws.Range("BK21").Value = 125.7
$ws.Range("B21").Copy()
$ws.Range("BK21").PasteSpecial(-4122)
$ws.Range("BK22").Value = 182.3
$ws.Range("A22").Copy()
$ws.Range("BK22").PasteSpecial(-4122)
$ws.Range("BK23").Value = 152.3
$ws.Range("A23").Copy()
$ws.Range("BK23").PasteSpecial(-4122)
$ws.Range("BK24").Value = 133
$ws.Range("D24").Copy()
$ws.Range("BK24").PasteSpecial(-4122)
$ws.Range("BK25").Value = 122.5
$ws.Range("E25").Copy()
$ws.Range("BK25").PasteSpecial(-4122)
$ws.Range("BK26").Value = 149.9
$ws.Range("A26").Copy()
$ws.Range("BK26").PasteSpecial(-4122)
$ws.Range("BK27").Value = 145.3
$ws.Range("A27").Copy()
$ws.Range("BK27").PasteSpecial(-4122)
$ws.Range("BK28").Value = 165.1
$ws.Range("A28").Copy()
$ws.Range("BK28").PasteSpecial(-4122)
$ws.Range("BK29").Value = 120.4
$ws.Range("H29").Copy()
$ws.Range("BK29").PasteSpecial(-4122)
$ws.Range("BK30").Value = 132.1
$ws.Range("E30").Copy()
$ws.Range("BK30").PasteSpecial(-4122)
$ws.Range("BK31").Value = 150.2
$ws.Range("A31").Copy()
$ws.Range("BK31").PasteSpecial(-4122)
$ws.Range("BK32").Value = 125.4
$ws.Range("H32").Copy()
$ws.Range("BK32").PasteSpecial(-4122)
$ws.Range("BK33").Value = 154
$ws.Range("A33").Copy()
$ws.Range("BK33").PasteSpecial(-4122)
$ws.Range("BK34").Value = 165.2
$ws.Range("A34").Copy()
$ws.Range("BK34").PasteSpecial(-4122)
$ws.Range("BK35").Value = 189
$ws.Range("A35").Copy()
$ws.Range("BK35").PasteSpecial(-4122)
$ws.Range("BK36").Value = 112.9
$ws.Range("F36").Copy()
$ws.Range("BK36").PasteSpecial(-4122)
$ws.Range("BK37").Value = 196.9
$ws.Range("A37").Copy()
$ws.Range("BK37").PasteSpecial(-4122)
$ws.Range("BK38").Value = 166.9
$ws.Range("A38").Copy()
$ws.Range("BK38").PasteSpecial(-4122)
$ws.Range("BK39").Value = 135.6
$ws.Range("D39").Copy()
$ws.Range("BK39").PasteSpecial(-4122)
$ws.Range("BK40").Value = 110.1
$ws.Range("J40").Copy()
$ws.Range("BK40").PasteSpecial(-4122)
$ws.Range("BK41").Value = 139.7
$ws.Range("C41").Copy()
$ws.Range("BK41").PasteSpecial(-4122)
$ws.Range("BK42").Value = 132.9
$ws.Range("K42").Copy()
$ws.Range("BK42").PasteSpecial(-4122)
$ws.Range("BK43").Value = 155.4
$ws.Range("A43").Copy()
$ws.Range("BK43").PasteSpecial(-4122)
$ws.Range("BK44").Value = 146.4
$ws.Range("A44").Copy()
$ws.Range("BK44").PasteSpecial(-4122)
$ws.Range("BK45").Value = 116
$ws.Range("AJ45").Copy()
$ws.Range("BK45").PasteSpecial(-4122)
$ws.Range("BK46").Value = 165.3
$ws.Range("A46").Copy()
$ws.Range("BK46").PasteSpecial(-4122)
$ws.Range("BK47").Value = 133.1
$ws.Range("B47").Copy()
$ws.Range("BK47").PasteSpecial(-4122)
$ws.Range("BK48").Value = 140.1
$ws.Range("A48").Copy()
$ws.Range("BK48").PasteSpecial(-4122)
$ws.Range("BK49").Value = 183.5
$ws.Range("A49").Copy()
$ws.Range("BK49").PasteSpecial(-4122)
$ws.Range("BK50").Value = 141.3
$ws.Range("A50").Copy()
$ws.Range("BK50").PasteSpecial(-4122)
$ws.Range("BK51").Value = 146.5
$ws.Range("A51").Copy()
$ws.Range("BK51").PasteSpecial(-4122)
$ws.Range("BK52").Value = 129.9
$ws.Range("C52").Copy()
$ws.Range("BK52").PasteSpecial(-4122)
$ws.Range("BK53").Value = 132.7
$ws.Range("C53").Copy()
$ws.Range("BK53").PasteSpecial(-4122)

$excel.CutCopyMode = 0